# Update the cryptos price/volume(1h) columns with the latest scraped values.
# Numeric-looking "Price" strings are prefixed with a leading apostrophe so
# Excel stores them as text (preserving exact formatting such as trailing
# zeros) instead of silently re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.788.44"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "2.413.46"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D5").Value = "'550.94"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'136.96"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'24.72"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "2.845.49"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "59.754.49"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "2.411.11"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D19").Value = "'4.35"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'329.53"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'65.76"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "'0.171"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").Value = "'8.58"
$ws.Range("E25").Value = "  +3.77%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "0.0₃0774"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "'1.76"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'170.62"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'18.58"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").Value = "'39.37"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'314.72"
$ws.Range("E40").Value = "  +8.98%  "
$ws.Range("D41").Value = "'0.408"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").Value = "'3.64"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "'137.92"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").Value = "'0.0963"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "'19.33"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").Value = "'0.408"
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "'17.50"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  -0.26%  "
